$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

$ws.Range("C6").Value = 1
$ws.Range("E6").Value = 15250

$ws.Range("F10").Select()
